$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.569.06'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.410.70'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.40'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.24'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('E10').Value = '  +2.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.50'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.781.87'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.425.78'
$ws.Range('E16').Value = '  +3.18%  '
$ws.Range('E17').Value = '  +3.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.577.56'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.42'
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.14'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0901'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.37'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.90'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.94'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.43'
$ws.Range('E29').Value = '  +3.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.43'
$ws.Range('E30').Value = '  +3.61%  '
$ws.Range('E31').Value = '  +16.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.39'
$ws.Range('E32').Value = '  +7.16%  '
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +3.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '131.47'
$ws.Range('E36').Value = '  +26.76%  '
$ws.Range('E37').Value = '  +3.07%  '
$ws.Range('E38').Value = '  +6.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.39'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.35'
$ws.Range('E42').Value = '  -4.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.945.11'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.83'
$ws.Range('E46').Value = '  +3.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.31'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.639.36'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('E49').Value = '  +4.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.69'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.20'
$ws.Range('E51').Value = '  -0.02%  '
